$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet between "2021-Q4" and "总计".
#    Worksheets.Add() always inserts at the front, so add then Move it into
#    place. Re-fetch sheet references by name after any operation that
#    changes sheet order/count -- handles here are positional, not stable.
# ---------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"
$newSheet.Move($wb.Worksheets.Item("总计"))
$ws = $wb.Worksheets.Item("2022-Q1")

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q1" with the same layout as "2021-Q4":
#    B:H header row, A2:A5 bold/boxed index column.
#    Pull the header + index-column formatting from "2021-Q4" via
#    copy/PasteSpecial(xlPasteFormats) so styles.xml isn't bloated with
#    near-duplicate entries.
# ---------------------------------------------------------------------------
$src.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

$src.Range("A2").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "162204"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "泰达宏利行业精选混合"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "7.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "75.56"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2.65"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0.1929"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 7

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "003501"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "泰达宏利睿智稳健灵活配置混合"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "73.79"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2.66"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0.0950"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 6

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "010181"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "兴业优势产业混合A"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.07"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "79.94"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "4.93"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0.0528"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 5

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "010182"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "兴业优势产业混合C"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "79.94"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "4.93"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0.0192"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 5

# ---------------------------------------------------------------------------
# 3. Update "总计": insert a "2022-Q1" row above the existing "2021-Q4" row
#    (new data on top, matching the diff) without using Rows.Insert(), which
#    pollutes styles.xml with an inherited-format duplicate. Shift the old
#    row 2 down to row 3 by hand, then overwrite row 2 with the new totals.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$oldB = $total.Range("B2").Value2
$oldC = $total.Range("C2").Value2
$oldD = $total.Range("D2").Value2

# Propagate A2's index-column style (bold/boxed) down to A3 first.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldB
$total.Range("C3").Value = $oldC
$total.Range("D3").Value = $oldD

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.36
